$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where column A holds the "as of" date that moved from 2025/11/23 -> 2025/11/24.
# The cell must stay a plain text value (it was stored as inline/shared text,
# not a real Excel date), so we force Text format before the write and then
# restore the cell's original (Normal) style so no formatting residue is left
# behind - otherwise Excel's automatic date recognition would silently turn
# the string into a date serial number.
$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value() -eq "2025/11/23") {
        $cell.NumberFormat = "@"
        $cell.Value = "2025/11/24"
        $cell.Style = "Normal"
    }
}
